$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FilesTab (row 4) query: reorder the RETURN columns (Format/Size moved earlier)
# to fix the case-files comparison.
$ws.Range("B4").Value = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN [''West Highland White Terrier''] 

OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '''') AS `File Name`,
       coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_type, '''') AS `File Type`,
       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(samp.sample_id, '''') AS `Sample ID`,
        coalesce(c.case_id, '''') AS `Case ID`,
        coalesce(demo.breed,'''') AS Breed ,
        coalesce(diag.disease_term,'''') AS Diagnosis'

# StudyFilesTab (row 5) query: unchanged text, now stored as a distinct shared
# string after the FilesTab query above was edited.
$ws.Range("B5").Value = 'MATCH (f:file)-->(s:study)<--(c:case)<--(demo:demographic)
WHERE demo.breed IN [''West Highland White Terrier'']
WITH
    f, s,
    [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, s,
    f.file_size /(1024^i) AS value, 10^precision AS factor,
    units[i] as unit
WITH
    f, s, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '''') AS `File Name`,
  coalesce(f.file_type, '''') AS `File Type`,
  coalesce("study", '''') AS `Association`,
  coalesce(f.file_description, '''') AS `Description`,
  coalesce(f.file_format, '''') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
  coalesce(s.clinical_study_designation,'''') AS `Study Code`'

# Update the window selection to match: active cell moves from B5 to B4.
$ws.Range("B4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
